# Final version of the report: fill in the last two journal entries
# (2024-07-18 and 2024-07-19, category "Analyse et état de l'art") and
# move the sheet's view/selection back to the top (E23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Row 22 - 2024-07-18, 10:00 -> 12:00
$ws.Range("A22").Value = 45491
$ws.Range("B22").Value = 0.41666666666666669
$ws.Range("C22").Value = 0.5
$ws.Range("E22").Value = "Analyse et état de l'art"

# Row 23 - 2024-07-19, 14:00 -> 16:00
$ws.Range("A23").Value = 45492
$ws.Range("B23").Value = 0.58333333333333337
$ws.Range("C23").Value = 0.66666666666666663
$ws.Range("E23").Value = "Analyse et état de l'art"

# Move the view back to the top and select the last entry typed in.
$ws.Application.Goto($ws.Range("A1"), $true)
[void]$ws.Range("E23").Select()
